# Refresh the "Price" / "Volume(1h)" columns (and, for two rows, the
# Coin/Link/Price text too) of the crypto-prices sheet to the latest
# scrape. All source values are plain text, not numbers (e.g. the
# thousands separator is "." and percentages keep their padding), so
# any value that Excel's automatic type detection would otherwise read
# as a genuine number (e.g. "594.40", "1.00", "0.0737") is written with
# a leading apostrophe, exactly as typing it into the grid by hand
# would, to keep it stored as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.830.45"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.496.21"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'594.40"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").Value = "'172.54"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.131"
$ws.Range("E9").Value = "  +4.37%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "4.100.05"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "'29.49"
$ws.Range("E14").Value = "  +5.19%  "
$ws.Range("D15").Value = "66.857.71"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "3.476.38"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'6.27"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'14.21"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "'394.13"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'73.32"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'6.13"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "'2.05"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'23.67"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").Value = "'162.23"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "'6.87"
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.837.27"
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0737"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").Value = "'26.14"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "'42.69"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").Value = "'337.40"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "'34.57"
$ws.Range("E48").Value = "  +2.58%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'6.42"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  -3.05%  "
